# Update Name of Algo
# Apply updated KNN imputation result values to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "C3"   = -13.228
    "B12"  = 4.976999999999999
    "C14"  = -11.895
    "C26"  = -13.04
    "C31"  = -12.68
    "B32"  = 6.95
    "C35"  = -12.311
    "B36"  = 8.673999999999999
    "C37"  = -13.387
    "B38"  = 5.595000000000001
    "C45"  = -12.67
    "B46"  = 6.247
    "B54"  = 5.340000000000001
    "B55"  = 4.742
    "C57"  = -13.697
    "B67"  = 5.34
    "B69"  = 5.002
    "B72"  = 5.380999999999999
    "B91"  = 6.378
    "B99"  = 5.718000000000001
    "C100" = -13.144
    "C102" = -12.916
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
